# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts: update the DAMSLTag (column I) and DialogAct
# (column J) values for the rows whose annotations changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 27;  Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 30;  Tag = "b";   Act = "Acknowledge (Backchannel)" },
    @{ Row = 49;  Tag = "qy";  Act = "Yes-No-Question" },
    @{ Row = 57;  Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 65;  Tag = "sv";  Act = "Statement-opinion" },
    @{ Row = 70;  Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 71;  Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 83;  Tag = "%";   Act = "Uninterpretable" },
    @{ Row = 99;  Tag = "sv";  Act = "Statement-opinion" },
    @{ Row = 103; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 109; Tag = "ba";  Act = "Appreciation" },
    @{ Row = 128; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 132; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 140; Tag = "sv";  Act = "Statement-opinion" },
    @{ Row = 153; Tag = "b";   Act = "Acknowledge (Backchannel)" },
    @{ Row = 158; Tag = "ba";  Act = "Appreciation" },
    @{ Row = 179; Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 184; Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 188; Tag = "ba";  Act = "Appreciation" },
    @{ Row = 207; Tag = "sv";  Act = "Statement-opinion" },
    @{ Row = 221; Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 224; Tag = "b";   Act = "Acknowledge (Backchannel)" },
    @{ Row = 231; Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 236; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 246; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 254; Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 263; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 292; Tag = "sv";  Act = "Statement-opinion" },
    @{ Row = 307; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 313; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 320; Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 356; Tag = "%";   Act = "Uninterpretable" },
    @{ Row = 363; Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 376; Tag = "sv";  Act = "Statement-opinion" },
    @{ Row = 379; Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 392; Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 394; Tag = "sv";  Act = "Statement-opinion" },
    @{ Row = 398; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 400; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 410; Tag = "sv";  Act = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
